$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target values (rows 1-20, columns B..I) after the "full search method" rework.
$data = @(
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 900, 900, 90, 130, "3month charter, 40 ships", 0.2283400748953568),
    @(20, 150, 800, 1000, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(30, 150, 700, 900, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(30, 150, 700, 900, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(30, 150, 700, 900, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(80, 130, 700, 1500, 80, 100, "NOT ADAPTED", -0.1364463133780686),
    @(30, 150, 700, 900, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(20, 150, 800, 1000, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(20, 150, 800, 1000, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(20, 150, 800, 1000, 90, 130, "NOT ADAPTED", -0.1364463133780686),
    @(20, 150, 800, 1000, 90, 130, "NOT ADAPTED", -0.1364463133780686),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 6).Value = $vals[4]   # F
    $ws.Cells.Item($row, 7).Value = $vals[5]   # G
    $ws.Cells.Item($row, 8).Value = $vals[6]   # H
    $ws.Cells.Item($row, 9).Value = $vals[7]   # I
}
